# Refresh crypto price (D) and 1h volume change % (E) columns with latest scraped values.
# Price values are forced to Text format before assignment (then formats are cleared again)
# so numeric-looking strings such as "573.67" or "61.158.19" are stored verbatim instead of
# being auto-converted into floating point numbers by Excel.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "61.158.19"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.52%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.399.64"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.54%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "573.67"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.29%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "139.07"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.17%  "
$ws.Range("E7").Value = "  +0.07%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.398.49"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.471"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.67%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "7.71"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +3.02%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.123"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  -1.77%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.381"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -2.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.976.40"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.50%  "
$ws.Range("E14").Value = "  -1.08%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "26.58"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.43%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0000172"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  -2.31%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.396.17"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +0.63%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "61.213.03"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.53%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.94"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.33%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "13.86"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  -1.41%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "9.32"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.53%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "375.26"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.530.63"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.56%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.553"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  -0.53%  "
$ws.Range("E25").Value = "  +0.15%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "71.09"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("E27").Value = "  -3.41%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.62"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -6.18%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.175"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +8.74%  "
$ws.Range("E30").Value = "  -0.54%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "7.42"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.54%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "8.12"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -2.01%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.15"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  -0.93%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "23.47"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -0.79%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.57"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.73%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "5.12"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.24%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.88"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.58%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "166.41"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.27%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.0771"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.47%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "25.94"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +5.43%  "
$ws.Range("E42").Value = "  +0.65%  "
$ws.Range("E43").Value = "  +0.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.777"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.23%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.94"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.69%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "4.40"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -0.21%  "
$ws.Range("E47").Value = "  -2.78%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "2.516.02"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +6.35%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "23.84"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +5.10%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "6.78"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  -0.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0261"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -0.10%  "
